# Update figures with new data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "117/230, 50.9% (95%CI 44.4-57.3%)"
$ws.Range("C3").Value = "5/230, 2.2% (95%CI 0.9-5%)"
$ws.Range("D3").Value = "172/230, 74.8% (95%CI 68.8-80%)"

$ws.Range("B7").Value = "1092/2112, 51.7% (95%CI 49.6-53.8%)"
$ws.Range("C7").Value = "42/2112, 2% (95%CI 1.5-2.7%)"
$ws.Range("D7").Value = "1638/2112, 77.6% (95%CI 75.7-79.3%)"

$ws.Range("B8").Value = "707/1437, 49.2% (95%CI 46.6-51.8%)"
$ws.Range("C8").Value = "10/1437, 0.7% (95%CI 0.4-1.3%)"
$ws.Range("D8").Value = "1087/1437, 75.6% (95%CI 73.4-77.8%)"

$ws.Range("B10").Value = "810/1586, 51.1% (95%CI 48.6-53.5%)"
$ws.Range("C10").Value = "32/1586, 2% (95%CI 1.4-2.8%)"
$ws.Range("D10").Value = "1210/1586, 76.3% (95%CI 74.1-78.3%)"

$ws.Range("B12").Value = "1018/1954, 52.1% (95%CI 49.9-54.3%)"
$ws.Range("C12").Value = "39/1954, 2% (95%CI 1.5-2.7%)"
$ws.Range("D12").Value = "1524/1954, 78% (95%CI 76.1-79.8%)"

$ws.Range("B13").Value = "333/704, 47.3% (95%CI 43.6-51%)"
$ws.Range("C13").Value = "22/704, 3.1% (95%CI 2.1-4.7%)"
$ws.Range("D13").Value = "515/704, 73.2% (95%CI 69.8-76.3%)"
